# Fruta / hortaliza, semanal
#
# New weekly price records arrived for "Naranja - Valencia - Primera"
# (Agrícola del Norte S.A. de Arica, Región de Coquimbo). Two identical
# records are inserted at the top of the existing data block (rows 59-60),
# pushing the previously existing rows 59-116 down to rows 61-118.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 59 (shifts old rows 59-116 down to 61-118).
$ws.Range("A59:A60").EntireRow.Insert()

# Populate the two new rows with the new price record.
$ws.Range("A59:A60").Value = 1
$ws.Range("B59:B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C59:C60").Value = "Arica y Parinacota"
$ws.Range("D59:D60").Value = 44923
$ws.Range("E59:E60").Value = 15
$ws.Range("F59:F60").Value = "Fruta"
$ws.Range("G59:G60").Value = 100102
$ws.Range("H59:H60").Value = "Cítricos"
$ws.Range("I59:I60").Value = 100102005
$ws.Range("J59:J60").Value = "Naranja"
$ws.Range("K59:K60").Value = "Valencia"
$ws.Range("L59:L60").Value = "Primera"
$ws.Range("M59:M60").Value = 350
$ws.Range("N59:N60").Value = 900
$ws.Range("O59:O60").Value = 1000
$ws.Range("P59:P60").Value = 943
$ws.Range("Q59:Q60").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R59:R60").Value = "Región de Coquimbo"
$ws.Range("S59:S60").Value = 943
$ws.Range("T59:T60").Value = 1
